$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) to reflect the new atom numbering scheme.
$ws.Range("C1").Value = "C3"
$ws.Range("D1").Value = "C4"
$ws.Range("E1").Value = "C5"
$ws.Range("F1").Value = "N1"
$ws.Range("G1").Value = "N2"
$ws.Range("H1").Value = "C6"

# Add two more mapped-atom columns (I, J) to the header, copying the
# formatting (bold, centered, bordered) from an existing header cell so
# the new cells share the same cell style as the rest of row 1.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$ws.Range("I1").Value = "C1"
$ws.Range("J1").Value = "C2"

# Populate the new data columns (I and J) for each data row, and keep
# the atom-mapping values in columns C-H consistent with the new header
# layout (the underlying mapped atoms are unchanged, only the header
# naming/columns shifted).
for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 3).Value = "C16"   # C
    $ws.Cells.Item($r, 4).Value = "C17"   # D
    $ws.Cells.Item($r, 5).Value = "C4"    # E
    $ws.Cells.Item($r, 6).Value = "N5"    # F
    $ws.Cells.Item($r, 7).Value = "N6"    # G
    $ws.Cells.Item($r, 8).Value = "C7"    # H
    $ws.Cells.Item($r, 9).Value = "C2"    # I
    $ws.Cells.Item($r, 10).Value = "C4"   # J
}
